# Add a new "Organisations" sheet to the PMHC upload workbook, after the
# existing "Practitioners" sheet, and populate it with the organisation
# upload-spec header row + one example/test row.

$wb = $excel.ActiveWorkbook

# --- Widen column A on the "Clients" sheet slightly ------------------------
$clients = $wb.Worksheets.Item("Clients")
$clients.Columns.Item(1).ColumnWidth = 14.333333333333334

# --- Create the new "Organisations" sheet after "Practitioners" ------------
$practitioners = $wb.Worksheets.Item("Practitioners")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $practitioners)
$ws.Name = "Organisations"

# Column widths (chosen to best match the authored widths under this
# runtime's 1/6-character rounding of ColumnWidth).
$ws.Columns.Item(1).ColumnWidth = 18.666666666666668
$ws.Columns.Item(2).ColumnWidth = 17.0
$ws.Columns.Item(3).ColumnWidth = 24.833333333333332
$ws.Columns.Item(4).ColumnWidth = 20.833333333333332
$ws.Columns.Item(5).ColumnWidth = 19.666666666666668

# Row 1: version marker
$ws.Range("A1").Value = "Version"
$ws.Range("B1").Value = 1

# Row 2: column headers
$ws.Range("A2").Value = "organisation_path"
$ws.Range("B2").Value = "organisation_key"
$ws.Range("C2").Value = "organisation_name"
$ws.Range("D2").Value = "organisation_legal_name"
$ws.Range("E2").Value = "organisation_abn"
$ws.Range("F2").Value = "organisation_type"
$ws.Range("G2").Value = "organisation_state"
$ws.Range("H2").Value = "organisation_status"
$ws.Range("I2").Value = "organisation_tags"

# Row 3: sample/test data
$ws.Range("A3").Value = "PHN999:NFP01"
$ws.Range("B3").Value = "NFP01"
$ws.Range("C3").Value = "Test Provider Organisation NFP1"
$ws.Range("E3").Value = 42072953425
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Page margins (inches -> points: 0.75in=54pt, 1in=72pt, 0.5in=36pt)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Match the authored selection/active cell on the new sheet
$ws.Range("F4").Select() | Out-Null
